# Applies the "STATUS GLOBAL" data refresh described by the commit:
# rows get re-sorted/re-paired within their PO groups, row 31/32 get replaced
# by the data that used to live in rows 33/34 (shifting the whole 33-37 block
# up by two), and the now-superfluous trailing rows 36-37 are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 2 / 3 (P-22/074) swap CER <-> DOS ---
$ws.Range("H2").Value = "22-074-DOS-0002"
$ws.Range("I2").Value = "FINAL QUALITY DOSSIER"
$ws.Range("J2").Value = "Dossier"
$ws.Range("P2").Value = " Enviado Rev.  //   Rev. "

$ws.Range("H3").Value = "22-074-CER-0002"
$ws.Range("I3").Value = "CERTIFICATES AND DECLARATIONS OF CONFORMITY"
$ws.Range("J3").Value = "Certificados"
$ws.Range("P3").Value = ""

# --- rows 4,5,6,8 (P-22/075) rotate: 4<-6, 6<-5, 5<-8, 8<-4 (row 7 untouched) ---
$ws.Range("H4").Value = "22-075-VDB-0002"
$ws.Range("I4").Value = "VENDOR DATA BOOK 22-075"
$ws.Range("J4").Value = "Dossier"
$ws.Range("K4").Value = "No"
$ws.Range("P4").Value = "   Rev. "

$ws.Range("H5").Value = "22-075-CER-0004"
$ws.Range("I5").Value = "TEST AND CERTIFICATION OF MATERIALS 22-075"
$ws.Range("J5").Value = "Certificados"

$ws.Range("H6").Value = "22-075-DOS-0002"
$ws.Range("I6").Value = "FINAL QUALITY DOSSIER 22-075"

$ws.Range("H8").Value = "22-075-MAN-0001"
$ws.Range("I8").Value = "INSTALLATION, OPERATION AND MAINTENANCE MANUAL 22-075"
$ws.Range("J8").Value = "Manual"
$ws.Range("K8").Value = "Sí"
$ws.Range("P8").Value = "   Rev.  // 22-04-2024 Enviado Rev.  //    Rev. "

# --- rows 10,12,13 (P-23/027) rotate: 10<-12, 12<-13, 13<-10 (row 11 untouched) ---
$ws.Range("H10").Value = "23-027-DOS-0002"
$ws.Range("I10").Value = "FINAL QUALITY DOSSIER"
$ws.Range("J10").Value = "Dossier"

$ws.Range("H12").Value = "23-027-LIS-0016"
$ws.Range("I12").Value = "SPARE PARTS LIST FOR PRECOMMISSIONING, COMMISSIONING AND START-UP"
$ws.Range("J12").Value = "Repuestos"

$ws.Range("H13").Value = "23-027-LIS-0017"
$ws.Range("I13").Value = "SPARE PARTS LIST FOR 2 YEARS OF OPERATION"

# --- rows 17 / 19 (P-23/037 level gauges) swap ---
$ws.Range("H17").Value = "23-037-LIS-0016"
$ws.Range("I17").Value = "NFXP3 - SPARE PARTS LIST FOR PRECOMMISSIONING, COMMISSIONING AND START-UP - LEVEL GAUGES"
$ws.Range("J17").Value = "Repuestos"

$ws.Range("H19").Value = "23-037-DOS-0002"
$ws.Range("I19").Value = "NFXP3 - MANUFACTURING RECORDS BOOK - LEVEL GAUGES"
$ws.Range("J19").Value = "Dossier"

# --- rows 24 / 25 (P-23/074 orifice plates) swap ---
$ws.Range("G24").Value = "5022_20-1043010910-00007"
$ws.Range("H24").Value = "23-074-LIS-0024"
$ws.Range("I24").Value = "NFXP4 - PRELIMINARY CARGO LIST FOR ORIFICE PLATES AND RESTRICTION ORIFICE"
$ws.Range("J24").Value = "Packing"

$ws.Range("G25").Value = "5022_20-1043010910-00004"
$ws.Range("H25").Value = "23-074-DOS-0002"
$ws.Range("I25").Value = "NFXP4 - MANUFACTURING RECORDS BOOK FOR ORIFICE PLATES AND RESTRICTION ORIFICE"
$ws.Range("J25").Value = "Dossier"

# --- rows 31-35: replaced by the former content of rows 33-37 (P-24/044 block) ---
# (dates/PO numbers are prefixed with a leading apostrophe so the COM layer
# keeps storing them as plain text instead of coercing to date/number types,
# matching the original "General"-formatted text cells)
$ws.Range("A31").Value = "P-24/044-S00"
$ws.Range("B31").Value = "'12-06-2024"
$ws.Range("C31").Value = "'12-09-2024"
$ws.Range("D31").Value = "'7080113517"
$ws.Range("E31").Value = "CEPSA"
$ws.Range("F31").Value = "Nivel"
$ws.Range("G31").Value = "V-MAN1642-2206-400-CER-001"
$ws.Range("H31").Value = "24-044-CER-0001"

$ws.Range("A32").Value = "P-24/044-S00"
$ws.Range("B32").Value = "'12-06-2024"
$ws.Range("C32").Value = "'12-09-2024"
$ws.Range("D32").Value = "'7080113517"
$ws.Range("E32").Value = "CEPSA"
$ws.Range("F32").Value = "Nivel"
$ws.Range("G32").Value = "V-MAN1642-2206-400-SPL-001"
$ws.Range("H32").Value = "24-044-SPL-0001"
$ws.Range("I32").Value = "LISTA DE REPUESTOS PARA 2 AÑOS"
$ws.Range("J32").Value = "Repuestos"

$ws.Range("G33").Value = "V-MAN1642-2206-400-LIS-001"
$ws.Range("H33").Value = "24-044-LIS-0002"
$ws.Range("I33").Value = "PACKING LIST Y AUTORIZACIÓN DE ENVÍO PARA FIRMA POR CEPSA"
$ws.Range("J33").Value = "Packing"

$ws.Range("G34").Value = "V-MAN1642-2206-400-DOS-001"
$ws.Range("H34").Value = "24-044-DOS-0001"
$ws.Range("I34").Value = "DOSSIER FINAL"
$ws.Range("J34").Value = "Dossier"

$ws.Range("G35").Value = "V-MAN1642-2206-400-CER-002"
$ws.Range("H35").Value = "24-044-CER-0002"
$ws.Range("I35").Value = "CERTIFICACIONES CE DE CONFORMIDAD"
$ws.Range("J35").Value = "Certificados"

# --- remove the now-duplicated trailing rows 36 and 37 ---
$ws.Rows.Item(37).Delete()
$ws.Rows.Item(36).Delete()
